$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 182, shifting existing rows 182-254 down to 183-255
$ws.Rows.Item(182).Insert()

# Populate the new row 182 with the new record's data
$ws.Cells.Item(182, 1).Value = 11
$ws.Cells.Item(182, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(182, 3).Value = "Bíobío"
$ws.Cells.Item(182, 4).Value = 45141
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = 100112032
$ws.Cells.Item(182, 7).Value = "Zapallo italiano"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 50
$ws.Cells.Item(182, 11).Value = 16000
$ws.Cells.Item(182, 12).Value = 16000
$ws.Cells.Item(182, 13).Value = 16000
$ws.Cells.Item(182, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(182, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(182, 16).Value = 320
$ws.Cells.Item(182, 17).Value = 50
$ws.Cells.Item(182, 18).Value = "Hortaliza"
